# daily_log.xlsx maintenance edit:
#  1) Clear a handful of stray "touched but empty" cells left over on the
#     2025-12-05 sheet (row 3: I3, J3, M3, P3, Q3, R3).
#  2) Append a brand-new daily sheet "2025-12-09" (after the last existing
#     sheet) using the same header/column layout as every other daily sheet,
#     with one test record.

$wb = $excel.ActiveWorkbook

# --- 1) Tidy up 2025-12-05 ---------------------------------------------
$ws1205 = $wb.Worksheets.Item("2025-12-05")
$ws1205.Range("I3").ClearContents()
$ws1205.Range("J3").ClearContents()
$ws1205.Range("M3").ClearContents()
$ws1205.Range("P3").ClearContents()
$ws1205.Range("Q3").ClearContents()
$ws1205.Range("R3").ClearContents()

# --- 2) Add the new "2025-12-09" sheet at the end -----------------------
$sheetCount = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($sheetCount)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "2025-12-09"

# Header row (same 18 columns used by every other daily-log sheet)
$headers = @("Registro ID", "Tipo Operación", "Contraparte", "Producto", "Peso Bruto (kg)", "Peso Tara (kg)", "Merma (kg)", "Peso Neto (kg)", "Precio x Kg", "Importe", "Chofer/Transporte", "Patente", "Incoterm", "Fecha Operacion", "Hora Ingreso", "Hora Salida", "Remito", "Observaciones")
for ($col = 1; $col -le $headers.Count; $col++) {
    $ws.Cells.Item(1, $col).Value = $headers[$col - 1]
}

# Column widths matching the other daily sheets
$colWidths = @(12, 15, 25, 25, 18, 18, 12, 18, 15, 18, 18, 15, 10, 15, 15, 15, 14, 30)
for ($col = 1; $col -le $colWidths.Count; $col++) {
    $ws.Columns.Item($col).ColumnWidth = $colWidths[$col - 1]
}

# Freeze the header row, like the rest of the workbook
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# One test record (row 2)
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "Compra"
$ws.Range("C2").Value = "Proveedor Test Timezone"
$ws.Range("D2").Value = "Producto Test"
$ws.Range("E2").Value = 1000
$ws.Range("F2").Value = 50
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 950
$ws.Range("I2").Value = 100
$ws.Range("J2").Value = 95000
$ws.Range("K2").Value = "Chofer Test"
$ws.Range("L2").Value = "ABC123"
$ws.Range("N2").Value = "09/12/25"
$ws.Range("O2").Value = "10:44"
$ws.Range("P2").Value = "10:44"
$ws.Range("R2").Value = "Prueba de zona horaria"

$ws.Range("A1").Select()
